$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds numeric-looking text (e.g. "26.904.28", "1.08")
# that must stay literal strings, not be re-interpreted as numbers by Excel.
# We briefly force text formatting while assigning, then restore each cells
# original style so no formatting side effects leak into the saved file.

$style_D2 = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.904.28"
$ws.Range("D2").Style = $style_D2
$ws.Range("E2").Value = "  -0.08%  "
$style_D3 = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.549.58"
$ws.Range("D3").Style = $style_D3
$ws.Range("E3").Value = "  -0.15%  "
$style_D5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.22"
$ws.Range("D5").Style = $style_D5
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("E6").Value = "  +0.72%  "
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("E8").Value = "  +2.26%  "
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("E12").Value = "  -0.16%  "
$style_D13 = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.550.75"
$ws.Range("D13").Style = $style_D13
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("E14").Value = "  +1.00%  "
$ws.Range("E15").Value = "  +0.83%  "
$style_D16 = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.911.58"
$ws.Range("D16").Style = $style_D16
$ws.Range("E16").Value = "  -0.09%  "
$style_D17 = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.66"
$ws.Range("D17").Style = $style_D17
$ws.Range("E17").Value = "  -0.09%  "
$style_D18 = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.40"
$ws.Range("D18").Style = $style_D18
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("E19").Value = "  +2.59%  "
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("E24").Value = "  -1.01%  "
$style_D25 = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.72"
$ws.Range("D25").Style = $style_D25
$ws.Range("E25").Value = "  +0.59%  "
$style_D26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.64"
$ws.Range("D26").Style = $style_D26
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  +0.63%  "
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("E30").Value = "  +1.95%  "
$style_D31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.08"
$ws.Range("D31").Style = $style_D31
$ws.Range("E31").Value = "  -0.71%  "
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("E33").Value = "  +4.82%  "
$style_D34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.410.96"
$ws.Range("D34").Style = $style_D34
$ws.Range("E34").Value = "  +2.57%  "
$style_D35 = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.59"
$ws.Range("D35").Style = $style_D35
$ws.Range("E35").Value = "  +2.64%  "
$style_D36 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.965"
$ws.Range("D36").Style = $style_D36
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("E38").Value = "  +0.99%  "
$style_D39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.526"
$ws.Range("D39").Style = $style_D39
$ws.Range("E39").Value = "  +0.91%  "
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("E42").Value = "  +3.63%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$style_D43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.996"
$ws.Range("D43").Style = $style_D43
$ws.Range("E43").Value = "  +0.72%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$style_D44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.28"
$ws.Range("D44").Style = $style_D44
$ws.Range("E44").Value = "  +1.25%  "
$style_D45 = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.52"
$ws.Range("D45").Style = $style_D45
$ws.Range("E45").Value = "  +1.32%  "
$ws.Range("E46").Value = "  +0.41%  "
$style_D47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.684.80"
$ws.Range("D47").Style = $style_D47
$ws.Range("E47").Value = "  -0.15%  "
$style_D48 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.32"
$ws.Range("D48").Style = $style_D48
$ws.Range("E48").Value = "  +1.24%  "
$ws.Range("E49").Value = "  +1.58%  "
$style_D50 = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0101"
$ws.Range("D50").Style = $style_D50
$ws.Range("E50").Value = "  +3.69%  "
